$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last row (row 19, Damian Lillard / Milwaukee Bucks) since the
# table shrinks from 18 data rows to 17 data rows (A1:C19 -> A1:C18).
$ws.Rows.Item(19).Delete()

# Overwrite rows 2-18 with the new table contents.
$data = @(
    @("Jared McCain", "PG,SG", "Philadelphia 76ers"),
    @("Cade Cunningham", "PG,SG", "Detroit Pistons"),
    @("Damian Lillard", "PG", "Milwaukee Bucks"),
    @("Royce O'Neale", "SF,PF", "Phoenix Suns"),
    @("Gradey Dick", "SG,SF", "Toronto Raptors"),
    @("Bam Adebayo", "C", "Miami Heat"),
    @("LaMelo Ball", "PG,SG", "Charlotte Hornets"),
    @("Brandon Miller", "SG,SF", "Charlotte Hornets"),
    @("Cam Thomas", "SG,SF", "Brooklyn Nets"),
    @("Yves Missi", "C", "New Orleans Pelicans"),
    @("Julius Randle", "PF", "Minnesota Timberwolves"),
    @("Cameron Johnson", "SF,PF", "Brooklyn Nets"),
    @("Brandon Ingram", "SG,SF,PF", "New Orleans Pelicans"),
    @("Derrick White", "PG,SG", "Boston Celtics"),
    @("Anthony Davis", "PF,C", "Los Angeles Lakers"),
    @("Isaiah Hartenstein", "C", "Oklahoma City Thunder"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
